$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "SCRIPT/G01P03A/um2404.ssb"
$ws.Range("B21").Value = 246
$ws.Range("C21").Value = " Good luck! And please\nbe careful!"
$ws.Range("D21").Value = " Удачи! И прошу, берегите себя!"
$ws.Range("E21").Value = " Ôäàœé! É ðñïšô, áåñåãéóå òåáÿ!"
